$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the first new "Actual Sample Time" row into the `series` block,
#     directly below the existing "Target Point" row (old row 94 -> new row 95).
$ws.Rows.Item(95).Insert()
$ws.Cells.Item(95, 1).Value = "series"
$ws.Cells.Item(95, 2).Value = "Actual Sample Time"
$ws.Cells.Item(95, 3).Value = "time"

# --- Insert the second new "Actual Sample Time" row into the `conc_time_values`
#     block, directly below its own "Target Point" row (old row 98 -> new row 100
#     once the row above has already been inserted).
$ws.Rows.Item(100).Insert()
$ws.Cells.Item(100, 1).Value = "conc_time_values"
$ws.Cells.Item(100, 2).Value = "Actual Sample Time"
$ws.Cells.Item(100, 3).Value = "time"

# --- Re-establish the AutoFilter over the now-larger A1:C109 range (inserting
#     rows inside a filtered range does not auto-grow the filter bounds).
$ws.AutoFilterMode = $false | Out-Null
$ws.Range("A1:C109").AutoFilter() | Out-Null

# --- Keep the hidden _FilterDatabase defined name in sync with the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$C`$109"
    }
}

# --- Restore the view state saved in the edited workbook: scrolled so row 89
#     is at the top, with C99 as the active selected cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 89
$win.ScrollColumn = 1
$ws.Range("C99").Select() | Out-Null
